$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.275.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.385.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "676.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.62"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.932.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.305.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.384.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.04"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.07%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "553.32"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.59"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.98"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.689.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.29"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0697"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.337"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0420"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.41%  "
